$wb = $excel.ActiveWorkbook

# --- DOC_SRC sheet: fiscal year-end date input -------------------------
$docSrc = $wb.Worksheets.Item("DOC_SRC")
$docSrc.Range("C5").Value = 40909

# --- ExecDB sheet: rename a column header, drop the extra blank column --
$execDb = $wb.Worksheets.Item("ExecDB")
$execDb.Range("T3").Value = "TimeVestRsValue"
$execDb.Columns("AA").Delete()

# --- Window / selection state -------------------------------------------
# Select the new input cell on DOC_SRC first ...
$docSrc.Range("C11").Select()

# ... then return to ExecDB (which remains the active/selected tab) and
# leave its selection/scroll position where the edit happened.
$execDb.Activate()
$excel.ActiveWindow.ScrollColumn = 14
$execDb.Range("T3").Select()

# Sheet-tab area was widened in the saved view.
$excel.ActiveWindow.TabRatio = 0.133
